$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "reduction_factor"
$ws.Range("C9").Value = 20
$ws.Range("D9").Value = 0.5
$ws.Range("E9").Value = 2

$ws.Range("Q8").Select()
